$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 5 data rows (rows 27-31, CRN crn_00025..crn_00029)
$ws.Range("A27:I31").EntireRow.Delete()

# New data values for rows 2-26, columns B:I
$arr = New-Object 'object[,]' 25,8
$arr[0,0] = 881.946489
$arr[0,1] = 242.501927
$arr[0,2] = -0.140574
$arr[0,3] = 1.86771
$arr[0,4] = 54.6875
$arr[0,5] = 2
$arr[0,6] = 70
$arr[0,7] = 56
$arr[1,0] = 853.563088
$arr[1,1] = 279.778026
$arr[1,2] = 0.149223
$arr[1,3] = 4.035843
$arr[1,4] = 51.5625
$arr[1,5] = 2
$arr[1,6] = 66
$arr[1,7] = 60
$arr[2,0] = 896.538497
$arr[2,1] = 250.993333
$arr[2,2] = 0.07372099999999999
$arr[2,3] = 2.402375
$arr[2,4] = 50
$arr[2,5] = 2
$arr[2,6] = 64
$arr[2,7] = 62
$arr[3,0] = 856.996006
$arr[3,1] = 272.878538
$arr[3,2] = -0.087912
$arr[3,3] = 3.494699
$arr[3,4] = 46.875
$arr[3,5] = 0
$arr[3,6] = 60
$arr[3,7] = 68
$arr[4,0] = 858.837881
$arr[4,1] = 230.82183
$arr[4,2] = -1.106368
$arr[4,3] = 1.011564
$arr[4,4] = 54.6875
$arr[4,5] = 1
$arr[4,6] = 70
$arr[4,7] = 57
$arr[5,0] = 913.830796
$arr[5,1] = 230.691923
$arr[5,2] = 0.04297
$arr[5,3] = 4.832913
$arr[5,4] = 48.4375
$arr[5,5] = 2
$arr[5,6] = 62
$arr[5,7] = 64
$arr[6,0] = 861.931649
$arr[6,1] = 241.958135
$arr[6,2] = 0.173264
$arr[6,3] = 5.955587
$arr[6,4] = 48.4375
$arr[6,5] = 2
$arr[6,6] = 62
$arr[6,7] = 64
$arr[7,0] = 856.319508
$arr[7,1] = 248.282937
$arr[7,2] = 0.08157300000000001
$arr[7,3] = 3.69036
$arr[7,4] = 53.125
$arr[7,5] = 0
$arr[7,6] = 68
$arr[7,7] = 60
$arr[8,0] = 912.701227
$arr[8,1] = 252.094699
$arr[8,2] = 0.033837
$arr[8,3] = 2.180909
$arr[8,4] = 54.6875
$arr[8,5] = 2
$arr[8,6] = 70
$arr[8,7] = 56
$arr[9,0] = 896.733764
$arr[9,1] = 239.507313
$arr[9,2] = -0.182095
$arr[9,3] = 2.141869
$arr[9,4] = 51.5625
$arr[9,5] = 2
$arr[9,6] = 66
$arr[9,7] = 60
$arr[10,0] = 902.100735
$arr[10,1] = 245.857253
$arr[10,2] = 0.043037
$arr[10,3] = 3.383029
$arr[10,4] = 54.6875
$arr[10,5] = 2
$arr[10,6] = 70
$arr[10,7] = 56
$arr[11,0] = 884.348984
$arr[11,1] = 233.42138
$arr[11,2] = 0.491044
$arr[11,3] = 5.890886
$arr[11,4] = 46.875
$arr[11,5] = 2
$arr[11,6] = 60
$arr[11,7] = 66
$arr[12,0] = 891.0269479999999
$arr[12,1] = 269.382221
$arr[12,2] = -0.266904
$arr[12,3] = 2.821687
$arr[12,4] = 46.875
$arr[12,5] = 1
$arr[12,6] = 60
$arr[12,7] = 67
$arr[13,0] = 861.334028
$arr[13,1] = 256.452531
$arr[13,2] = 0.120206
$arr[13,3] = 2.916221
$arr[13,4] = 50
$arr[13,5] = 1
$arr[13,6] = 64
$arr[13,7] = 63
$arr[14,0] = 873.613148
$arr[14,1] = 231.483972
$arr[14,2] = 0.269043
$arr[14,3] = 5.332094
$arr[14,4] = 54.6875
$arr[14,5] = 1
$arr[14,6] = 70
$arr[14,7] = 57
$arr[15,0] = 895.562584
$arr[15,1] = 219.485318
$arr[15,2] = -1.399416
$arr[15,3] = 1.298312
$arr[15,4] = 50
$arr[15,5] = 2
$arr[15,6] = 64
$arr[15,7] = 62
$arr[16,0] = 866.6680249999999
$arr[16,1] = 263.282676
$arr[16,2] = 1.151636
$arr[16,3] = 5.725297
$arr[16,4] = 48.4375
$arr[16,5] = 0
$arr[16,6] = 62
$arr[16,7] = 66
$arr[17,0] = 856.671243
$arr[17,1] = 247.115408
$arr[17,2] = -0.533533
$arr[17,3] = 1.937697
$arr[17,4] = 53.125
$arr[17,5] = 1
$arr[17,6] = 68
$arr[17,7] = 59
$arr[18,0] = 874.655362
$arr[18,1] = 245.966022
$arr[18,2] = 0.326655
$arr[18,3] = 4.382966
$arr[18,4] = 48.4375
$arr[18,5] = 2
$arr[18,6] = 62
$arr[18,7] = 64
$arr[19,0] = 888.2242220000001
$arr[19,1] = 245.693509
$arr[19,2] = 0.188066
$arr[19,3] = 3.193398
$arr[19,4] = 51.5625
$arr[19,5] = 2
$arr[19,6] = 66
$arr[19,7] = 60
$arr[20,0] = 893.807727
$arr[20,1] = 254.932305
$arr[20,2] = 0.482877
$arr[20,3] = 5.395766
$arr[20,4] = 46.875
$arr[20,5] = 1
$arr[20,6] = 60
$arr[20,7] = 67
$arr[21,0] = 843.694317
$arr[21,1] = 252.224038
$arr[21,2] = 0.142732
$arr[21,3] = 7.942065
$arr[21,4] = 48.4375
$arr[21,5] = 1
$arr[21,6] = 62
$arr[21,7] = 65
$arr[22,0] = 886.506652
$arr[22,1] = 272.621458
$arr[22,2] = -0.449558
$arr[22,3] = 1.864205
$arr[22,4] = 54.6875
$arr[22,5] = 0
$arr[22,6] = 70
$arr[22,7] = 58
$arr[23,0] = 898.185211
$arr[23,1] = 257.13574
$arr[23,2] = -0.46345
$arr[23,3] = 4.197988
$arr[23,4] = 51.5625
$arr[23,5] = 2
$arr[23,6] = 66
$arr[23,7] = 60
$arr[24,0] = 869.643886
$arr[24,1] = 262.298498
$arr[24,2] = -0.933485
$arr[24,3] = 2.552094
$arr[24,4] = 45.3125
$arr[24,5] = 0
$arr[24,6] = 58
$arr[24,7] = 70

$ws.Range("B2:I26").Value = $arr
